$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.14538034605772
$ws.Range("C2").Value = 11.09381386923713
$ws.Range("D2").Value = 4.991507198785488
$ws.Range("E2").Value = 12.33569114822233
$ws.Range("F2").Value = 24.49992230631182
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 22.19670225903976
$ws.Range("L2").Value = 10.03587379587379
$ws.Range("M2").Value = 14.53851040075903
$ws.Range("O2").Value = 21.94641149117582
$ws.Range("B3").Value = 14.54781237916841
$ws.Range("C3").Value = 10.8148330778737
$ws.Range("D3").Value = 4.946680220603967
$ws.Range("E3").Value = 12.38353142493903
$ws.Range("F3").Value = 24.5541968581864
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 22.34220659126368
$ws.Range("L3").Value = 10.04495477195722
$ws.Range("M3").Value = 14.40417993689605
$ws.Range("O3").Value = 22.04265431265354
$ws.Range("B4").Value = 14.16874570717532
$ws.Range("C4").Value = 10.63888572589536
$ws.Range("D4").Value = 4.918746072883493
$ws.Range("E4").Value = 12.41445368153827
$ws.Range("F4").Value = 24.59586488378161
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 22.43698940105817
$ws.Range("L4").Value = 10.05196209405625
$ws.Range("M4").Value = 14.32267194227342
$ws.Range("O4").Value = 22.10802811107565
$ws.Range("B5").Value = 14.01142594475998
$ws.Range("C5").Value = 10.5660897946672
$ws.Range("D5").Value = 4.907265396023751
$ws.Range("E5").Value = 12.4274452087177
$ws.Range("F5").Value = 24.61493434888935
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 22.47698246070575
$ws.Range("L5").Value = 10.05517801652483
$ws.Range("M5").Value = 14.28972932047401
$ws.Range("O5").Value = 22.13624164720362
$ws.Range("B6").Value = 13.9851377957467
$ws.Range("C6").Value = 10.55393814562083
$ws.Range("D6").Value = 4.905353334729021
$ws.Range("E6").Value = 12.42962605882838
$ws.Range("F6").Value = 24.61822674386154
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 22.48370591301038
$ws.Range("L6").Value = 10.05573379460261
$ws.Range("M6").Value = 14.28427650101861
$ws.Range("O6").Value = 22.14102132014271
$ws.Range("B7").Value = 14.16663526937803
$ws.Range("C7").Value = 10.63790831086798
$ws.Range("D7").Value = 4.918591626597696
$ws.Range("E7").Value = 12.41462730732292
$ws.Range("F7").Value = 24.59611361342542
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 22.43752322159789
$ws.Range("L7").Value = 10.05200400545701
$ws.Range("M7").Value = 14.32222652604091
$ws.Range("O7").Value = 22.1084022470586
$ws.Range("B8").Value = 14.94198064445245
$ws.Range("C8").Value = 10.99862959271999
$ws.Range("D8").Value = 4.976139494801965
$ws.Range("E8").Value = 12.35186586394249
$ws.Range("F8").Value = 24.51689968799386
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 22.2457423653412
$ws.Range("L8").Value = 10.03870805523297
$ws.Range("M8").Value = 14.49200872306552
$ws.Range("O8").Value = 21.9782889380086
$ws.Range("B9").Value = 16.35795579589243
$ws.Range("C9").Value = 11.66610024122545
$ws.Range("D9").Value = 5.08548238843701
$ws.Range("E9").Value = 12.24102123337161
$ws.Range("F9").Value = 24.42809738405555
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 21.91285814054433
$ws.Range("L9").Value = 10.02397404126998
$ws.Range("M9").Value = 14.83140026703192
$ws.Range("O9").Value = 21.77323334470488
$ws.Range("B10").Value = 17.32542332544563
$ws.Range("C10").Value = 12.1285160297046
$ws.Range("D10").Value = 5.163344934550947
$ws.Range("E10").Value = 12.16696468596374
$ws.Range("F10").Value = 24.40380503533923
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 21.6946417625676
$ws.Range("L10").Value = 10.02003269903116
$ws.Range("M10").Value = 15.08298202964178
$ws.Range("O10").Value = 21.65347172931432
$ws.Range("B11").Value = 17.74819780227866
$ws.Range("C11").Value = 12.3321343486555
$ws.Range("D11").Value = 5.19816100724583
$ws.Range("E11").Value = 12.13486120757244
$ws.Range("F11").Value = 24.40170436372922
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 21.60109831678497
$ws.Range("L11").Value = 10.01972676734761
$ws.Range("M11").Value = 15.19755475806762
$ws.Range("O11").Value = 21.60576921201117
$ws.Range("B12").Value = 17.9056946738426
$ws.Range("C12").Value = 12.40822085606678
$ws.Range("D12").Value = 5.21125240203619
$ws.Range("E12").Value = 12.12293120608051
$ws.Range("F12").Value = 24.40219862726455
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 21.56649994658978
$ws.Range("L12").Value = 10.01982394753606
$ws.Range("M12").Value = 15.24092842233462
$ws.Range("O12").Value = 21.58868577796234
$ws.Range("B13").Value = 17.87189201278599
$ws.Range("C13").Value = 12.39188035618807
$ws.Range("D13").Value = 5.208437150275275
$ws.Range("E13").Value = 12.12549047124365
$ws.Range("F13").Value = 24.40203478763314
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 21.57391463169302
$ws.Range("L13").Value = 10.01979355709583
$ws.Range("M13").Value = 15.23158817200971
$ws.Range("O13").Value = 21.59232130345592
$ws.Range("B14").Value = 17.76120780072257
$ws.Range("C14").Value = 12.33841473345322
$ws.Range("D14").Value = 5.199239934953979
$ws.Range("E14").Value = 12.13387517803357
$ws.Range("F14").Value = 24.40171916775633
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 21.5982353555553
$ws.Range("L14").Value = 10.01973049785767
$ws.Range("M14").Value = 15.20112356201831
$ws.Range("O14").Value = 21.60434406471011
$ws.Range("B15").Value = 17.69306920776168
$ws.Range("C15").Value = 12.3055313274298
$ws.Range("D15").Value = 5.193594131999346
$ws.Range("E15").Value = 12.13904056929299
$ws.Range("F15").Value = 24.40169386266017
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 21.61323991821904
$ws.Range("L15").Value = 10.01971958989026
$ws.Range("M15").Value = 15.18246056745505
$ws.Range("O15").Value = 21.61183620975237
$ws.Range("B16").Value = 17.29743581572846
$ws.Range("C16").Value = 12.1150690998224
$ws.Range("D16").Value = 5.161056966076902
$ws.Range("E16").Value = 12.1690945313029
$ws.Range("F16").Value = 24.40412268095933
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 21.70087034708149
$ws.Range("L16").Value = 10.02008254825916
$ws.Range("M16").Value = 15.07549434672857
$ws.Range("O16").Value = 21.65672609030307
$ws.Range("B17").Value = 17.05020785426936
$ws.Range("C17").Value = 11.99646373021633
$ws.Range("D17").Value = 5.140937574740343
$ws.Range("E17").Value = 12.1879369162707
$ws.Range("F17").Value = 24.40790736741342
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 21.75609582839699
$ws.Range("L17").Value = 10.0206855263481
$ws.Range("M17").Value = 15.00988453547246
$ws.Range("O17").Value = 21.68600473278995
$ws.Range("B18").Value = 16.90638286997168
$ws.Range("C18").Value = 11.92761366493049
$ws.Range("D18").Value = 5.129308924387653
$ws.Range("E18").Value = 12.1989238288203
$ws.Range("F18").Value = 24.4109265579711
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 21.7883987716746
$ws.Range("L18").Value = 10.0211723603682
$ws.Range("M18").Value = 14.97216127987241
$ws.Range("O18").Value = 21.70348264058395
$ws.Range("B19").Value = 16.85741043344877
$ws.Range("C19").Value = 11.90419539091037
$ws.Range("D19").Value = 5.125362137911567
$ws.Range("E19").Value = 12.20266947862082
$ws.Range("F19").Value = 24.41209337161294
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 21.7994284774094
$ws.Range("L19").Value = 10.02136126506282
$ws.Range("M19").Value = 14.95939214259246
$ws.Range("O19").Value = 21.70950968738914
$ws.Range("B20").Value = 17.07669479760094
$ws.Range("C20").Value = 12.00915520226178
$ws.Range("D20").Value = 5.143085209449022
$ws.Range("E20").Value = 12.18591567095636
$ws.Range("F20").Value = 24.40741727811221
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 21.75016121801041
$ws.Range("L20").Value = 10.0206068526003
$ws.Range("M20").Value = 15.01686761057638
$ws.Range("O20").Value = 21.68282193821545
$ws.Range("B21").Value = 17.79378973575104
$ws.Range("C21").Value = 12.35414692718697
$ws.Range("D21").Value = 5.201943941711005
$ws.Range("E21").Value = 12.13140623712744
$ws.Range("F21").Value = 24.40177685447376
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 21.59106938119373
$ws.Range("L21").Value = 10.01974324509759
$ws.Range("M21").Value = 15.21007234049048
$ws.Range("O21").Value = 21.60078603295706
$ws.Range("B22").Value = 18.24726689532432
$ws.Range("C22").Value = 12.57365729951368
$ws.Range("D22").Value = 5.239868604478138
$ws.Range("E22").Value = 12.09710316684506
$ws.Range("F22").Value = 24.4056084467138
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 21.49190022235803
$ws.Range("L22").Value = 10.02042009936713
$ws.Range("M22").Value = 15.33625568715741
$ws.Range("O22").Value = 21.55288833282443
$ws.Range("B23").Value = 18.00665839134712
$ws.Range("C23").Value = 12.45706149036303
$ws.Range("D23").Value = 5.219679071771999
$ws.Range("E23").Value = 12.11529074733454
$ws.Range("F23").Value = 24.40287500108372
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 21.54438838692641
$ws.Range("L23").Value = 10.0199455586241
$ws.Range("M23").Value = 15.26892730825579
$ws.Range("O23").Value = 21.57792720897541
$ws.Range("B24").Value = 17.06472530731568
$ws.Range("C24").Value = 12.00341944669213
$ws.Range("D24").Value = 5.142114455665278
$ws.Range("E24").Value = 12.18682899564037
$ws.Range("F24").Value = 24.40763622079619
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 21.75284253262516
$ws.Range("L24").Value = 10.02064198431517
$ws.Range("M24").Value = 15.01371057024496
$ws.Range("O24").Value = 21.6842588702224
$ws.Range("B25").Value = 15.98710838541653
$ws.Range("C25").Value = 11.49019986289804
$ws.Range("D25").Value = 5.056312504876446
$ws.Range("E25").Value = 12.26970607836592
$ws.Range("F25").Value = 24.4449528860015
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 21.99828599568478
$ws.Range("L25").Value = 10.0267487205731
$ws.Range("M25").Value = 14.73907054651484
$ws.Range("O25").Value = 21.82330608012984
